$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 134.57143
$ws.Range("I11").Value = 134.57143
$ws.Range("K11").Value = 134.57143
$ws.Range("M11").Value = 5.428570000000008
$ws.Range("H92").Value = 58878.207
$ws.Range("I92").Value = 33653.535
$ws.Range("K92").Value = 33653.535
$ws.Range("M92").Value = -32405.535
$ws.Range("H100").Value = 3365.85
$ws.Range("J100").Value = 3133.8572
$ws.Range("L100").Value = 3133.8572
$ws.Range("N100").Value = -4215.8572
$ws.Range("H106").Value = 5790.3125
$ws.Range("I106").Value = 4760.5835
$ws.Range("K106").Value = 4760.5835
$ws.Range("M106").Value = -4129.5835
$ws.Range("H132").Value = 2209.4285
$ws.Range("I132").Value = 2019.4147
$ws.Range("K132").Value = 6058.2441
$ws.Range("M132").Value = -3528.2441
$ws.Range("H135").Value = 26266.742
$ws.Range("I135").Value = 31669.69
$ws.Range("K135").Value = 285027.21
$ws.Range("M135").Value = -282492.21
$ws.Range("H141").Value = 400.9365
$ws.Range("I141").Value = 400.9365
$ws.Range("K141").Value = 1202.8095
$ws.Range("M141").Value = 3977.1905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 21612.385
$ws.Range("I61").Value = 3241.913
$ws.Range("J61").Value = 36182.07
$ws.Range("K61").Value = 3241.913
$ws.Range("L61").Value = 36182.07
$ws.Range("M61").Value = -3029.913
$ws.Range("N61").Value = -36606.07
$ws.Range("H74").Value = 53093.63
$ws.Range("I74").Value = 28184.617
$ws.Range("K74").Value = 28184.617
$ws.Range("M74").Value = -27310.617
$ws.Range("H76").Value = 148284.3
$ws.Range("J76").Value = 148284.3
$ws.Range("L76").Value = 148284.3
$ws.Range("N76").Value = -148960.3
$ws.Range("H77").Value = 53093.63
$ws.Range("I77").Value = 28184.617
$ws.Range("K77").Value = 140923.085
$ws.Range("M77").Value = -136555.085
$ws.Range("H79").Value = 148284.3
$ws.Range("J79").Value = 148284.3
$ws.Range("L79").Value = 148284.3
$ws.Range("N79").Value = -150624.3
$ws.Range("H80").Value = 47500
$ws.Range("I80").Value = 40000
$ws.Range("J80").Value = 55000
$ws.Range("K80").Value = 40000
$ws.Range("L80").Value = 55000
$ws.Range("M80").Value = -39002
$ws.Range("N80").Value = -56996
$ws.Range("H83").Value = 47500
$ws.Range("I83").Value = 40000
$ws.Range("J83").Value = 55000
$ws.Range("K83").Value = 120000
$ws.Range("L83").Value = 165000
$ws.Range("M83").Value = -115008
$ws.Range("N83").Value = -174984
$ws.Range("H97").Value = 1056.08
$ws.Range("I97").Value = 614.6
$ws.Range("K97").Value = 614.6
$ws.Range("M97").Value = -118.6
$ws.Range("H122").Value = 73711.84
$ws.Range("I122").Value = 4781.875
$ws.Range("K122").Value = 14345.625
$ws.Range("M122").Value = -11895.625
$ws.Range("H136").Value = 21612.385
$ws.Range("I136").Value = 3241.913
$ws.Range("J136").Value = 36182.07
$ws.Range("K136").Value = 9725.739
$ws.Range("L136").Value = 108546.21
$ws.Range("M136").Value = -7175.739
$ws.Range("N136").Value = -113646.21

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 36597.96
$ws.Range("I94").Value = 557.1579
$ws.Range("J94").Value = 150727.17
$ws.Range("K94").Value = 557.1579
$ws.Range("L94").Value = 150727.17
$ws.Range("M94").Value = -106.1579
$ws.Range("N94").Value = -151629.17
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("H134").Value = 2439.8594
$ws.Range("I134").Value = 1933.1731
$ws.Range("K134").Value = 5799.5193
$ws.Range("M134").Value = -3264.5193
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 29500
$ws.Range("J53").Value = 34000
$ws.Range("L53").Value = 34000
$ws.Range("N53").Value = -35214
$ws.Range("H58").Value = 7899
$ws.Range("J58").Value = 25162.6
$ws.Range("L58").Value = 25162.6
$ws.Range("N58").Value = -25568.6
$ws.Range("H105").Value = 1692.0952
$ws.Range("I105").Value = 1362.3334
$ws.Range("K105").Value = 1362.3334
$ws.Range("M105").Value = 384.6666
$ws.Range("H107").Value = 40038180
$ws.Range("I107").Value = 52679940
$ws.Range("J107").Value = 5942.6665
$ws.Range("K107").Value = 52679940
$ws.Range("L107").Value = 5942.6665
$ws.Range("M107").Value = -52678020
$ws.Range("N107").Value = -9782.666499999999
$ws.Range("H136").Value = 7899
$ws.Range("J136").Value = 25162.6
$ws.Range("L136").Value = 75487.79999999999
$ws.Range("N136").Value = -80587.79999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 11906.728
$ws.Range("I3").Value = 8046.5293
$ws.Range("J3").Value = 25031.4
$ws.Range("K3").Value = 24139.5879
$ws.Range("L3").Value = 75094.20000000001
$ws.Range("M3").Value = -24027.5879
$ws.Range("N3").Value = -75318.20000000001
$ws.Range("H28").Value = 676.4
$ws.Range("I28").Value = 470.75
$ws.Range("K28").Value = 1412.25
$ws.Range("M28").Value = -1180.25
$ws.Range("H113").Value = 1440.762
$ws.Range("I113").Value = 683.4286
$ws.Range("J113").Value = 1819.4286
$ws.Range("K113").Value = 2050.2858
$ws.Range("L113").Value = 5458.2858
$ws.Range("M113").Value = 119.7142000000003
$ws.Range("N113").Value = -9798.2858
$ws.Range("H131").Value = 2785
$ws.Range("I131").Value = 1402.0769
$ws.Range("J131").Value = 4419.364
$ws.Range("K131").Value = 4206.2307
$ws.Range("L131").Value = 13258.092
$ws.Range("M131").Value = 833.7692999999999
$ws.Range("N131").Value = -23338.092
$ws.Range("H139").Value = 2536.739
$ws.Range("I139").Value = 1492.85
$ws.Range("K139").Value = 4478.549999999999
$ws.Range("M139").Value = 661.4500000000007

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 71975.234
$ws.Range("I97").Value = 102640.664
$ws.Range("J97").Value = 2978
$ws.Range("K97").Value = 102640.664
$ws.Range("L97").Value = 2978
$ws.Range("M97").Value = -102144.664
$ws.Range("N97").Value = -3970
$ws.Range("H99").Value = 8226.25
$ws.Range("I99").Value = 2511.8
$ws.Range("J99").Value = 36798.5
$ws.Range("K99").Value = 2511.8
$ws.Range("L99").Value = 36798.5
$ws.Range("M99").Value = -265.8000000000002
$ws.Range("N99").Value = -41290.5
$ws.Range("H108").Value = 90684
$ws.Range("J108").Value = 90684
$ws.Range("L108").Value = 90684
$ws.Range("N108").Value = -98364
$ws.Range("H126").Value = 229005.25
$ws.Range("I126").Value = 229005.25
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 687015.75
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -684545.75
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 40234.582
$ws.Range("I40").Value = 55614.5
$ws.Range("J40").Value = 9474.75
$ws.Range("K40").Value = 55614.5
$ws.Range("L40").Value = 9474.75
$ws.Range("M40").Value = -55478.5
$ws.Range("N40").Value = -9746.75
$ws.Range("H82").Value = 1808.44
$ws.Range("I82").Value = 2031.5
$ws.Range("J82").Value = 1234.8572
$ws.Range("K82").Value = 2031.5
$ws.Range("L82").Value = 1234.8572
$ws.Range("M82").Value = -1670.5
$ws.Range("N82").Value = -1956.8572
$ws.Range("H85").Value = 1808.44
$ws.Range("I85").Value = 2031.5
$ws.Range("J85").Value = 1234.8572
$ws.Range("K85").Value = 2031.5
$ws.Range("L85").Value = 1234.8572
$ws.Range("M85").Value = -783.5
$ws.Range("N85").Value = -3730.8572
$ws.Range("H93").Value = 1670.4615
$ws.Range("I93").Value = 1476
$ws.Range("K93").Value = 1476
$ws.Range("M93").Value = -228
$ws.Range("H122").Value = 3504.375
$ws.Range("I122").Value = 3338.0667
$ws.Range("K122").Value = 10014.2001
$ws.Range("M122").Value = -7564.2001
$ws.Range("H132").Value = 3825.4443
$ws.Range("I132").Value = 2420.75
$ws.Range("K132").Value = 7262.25
$ws.Range("M132").Value = -4732.25
$ws.Range("H136").Value = 86235.55
$ws.Range("I136").Value = 3534.6667
$ws.Range("K136").Value = 10604.0001
$ws.Range("M136").Value = -8054.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1246.5
$ws.Range("I107").Value = 1246.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3739.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1819.5
$ws.Range("H122").Value = 3468.5334
$ws.Range("I122").Value = 3141.7837
$ws.Range("K122").Value = 9425.3511
$ws.Range("M122").Value = -6975.3511
$ws.Range("H126").Value = 3309.889
$ws.Range("I126").Value = 3379.875
$ws.Range("K126").Value = 10139.625
$ws.Range("M126").Value = -7669.625
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H136").Value = 6038.614
$ws.Range("I136").Value = 6270.3335
$ws.Range("J136").Value = 4995.875
$ws.Range("K136").Value = 18811.0005
$ws.Range("L136").Value = 14987.625
$ws.Range("M136").Value = -16261.0005
$ws.Range("N136").Value = -20087.625
$ws.Range("N107").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
